$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The backend re-synced the "New Books Report" sheet: every row got a fresh
# id and a refreshed createdAt/updatedAt timestamp, the rows were
# re-ordered, the "Fundamentals of Wavelets" row was soft-deleted
# (isDeleted = TRUE, deletedAt populated), and a new row (How to Think Like
# Sherlock Holmes) was appended. Rewrite rows 2-11 in full to match.

$ws.Range("A2").Value = "7e424a44-01ea-4d92-83fe-d7f003c8bde8"
$ws.Range("B2").Value = "Fundamentals of Wavelets"
$ws.Range("C2").Value = "signal_processing"
$ws.Range("D2").Value = "Wiley"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = 45391.35505844907
$ws.Range("G2").Value = 45391.35505844907
$ws.Range("H2").Value = 45391.592312939814
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$ws.Range("A3").Value = "52d4f340-5f7e-4d86-8ab7-25566eb64995"
$ws.Range("B3").Value = "Data Smart"
$ws.Range("C3").Value = "data_science"
$ws.Range("D3").Value = "Wiley"
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = 45391.35505844907
$ws.Range("G3").Value = 45391.35505844907

$ws.Range("A4").Value = "ef0fbff6-2d51-4a9b-9a34-212512cf8d4b"
$ws.Range("B4").Value = "God Created the Integers"
$ws.Range("C4").Value = "mathematics"
$ws.Range("D4").Value = "Penguin"
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 45391.35505844907
$ws.Range("G4").Value = 45391.35505844907

$ws.Range("A5").Value = "0eac702b-4e5b-48e8-9948-007b685d0915"
$ws.Range("B5").Value = "Superfreakonomics"
$ws.Range("C5").Value = "economics"
$ws.Range("D5").Value = "HarperCollins"
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = 45391.35505844907
$ws.Range("G5").Value = 45391.35505844907

$ws.Range("A6").Value = "d915b0f9-7293-4011-8b9a-2171d8984dd0"
$ws.Range("B6").Value = "Orientalism"
$ws.Range("C6").Value = "history"
$ws.Range("D6").Value = "Penguin"
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = 45391.35505844907
$ws.Range("G6").Value = 45391.35505844907

$ws.Range("A7").Value = "2e1b1f40-cb5b-434e-a609-d33f6c7074b3"
$ws.Range("B7").Value = "Nature of Statistical Learning Theory, The"
$ws.Range("C7").Value = "data_science"
$ws.Range("D7").Value = "Springer"
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 45391.35505844907
$ws.Range("G7").Value = 45391.35505844907

$ws.Range("A8").Value = "6f994baf-5553-4131-a175-73672ce6081f"
$ws.Range("B8").Value = "Integration of the Indian States"
$ws.Range("C8").Value = "history"
$ws.Range("D8").Value = "Orient Blackswan"
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = 45391.35505844907
$ws.Range("G8").Value = 45391.35505844907

$ws.Range("A9").Value = "0fe04d0b-4130-4b7f-a562-e540d252194b"
$ws.Range("B9").Value = "Drunkard's Walk, The"
$ws.Range("C9").Value = "science"
$ws.Range("D9").Value = "Penguin"
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = 45391.35505844907
$ws.Range("G9").Value = 45391.35505844907

$ws.Range("A10").Value = "0a0ee015-8983-47f8-b20a-f2ef4137db8a"
$ws.Range("B10").Value = "Image Processing & Mathematical Morphology"
$ws.Range("C10").Value = "signal_processing"
$ws.Range("D10").Value = "CRC"
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = 45391.35505844907
$ws.Range("G10").Value = 45391.35505844907

$ws.Range("A11").Value = "3894fa3a-6ca0-4c38-9e49-1f0fa3cdfbb3"
$ws.Range("B11").Value = "How to Think Like Sherlock Holmes"
$ws.Range("C11").Value = "psychology"
$ws.Range("D11").Value = "Penguin"
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = 45391.35505844907
$ws.Range("G11").Value = 45391.35505844907

# Row 11 is brand new, so F11/G11 need the same date number format the
# other createdAt/updatedAt cells use.
$ws.Range("F10:G10").Copy()
$ws.Range("F11:G11").PasteSpecial(-4122)

